# feat: upload to cloudinary
#
# Add an "image_url" column (F) to the siswa (student) export template so
# the generated spreadsheet has a place for the Cloudinary-hosted photo
# URL, alongside notelp / nama / nis / nisn / rombel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in F1, using the same text/fill/alignment formatting as the
# other header cells (copy the format from the neighbouring E1 header).
$ws.Range("F1").Value = "image_url"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Widen the new column so the image URL is readable.
$ws.Columns.Item(6).ColumnWidth = 16.7109375

# The new header cell becomes the active selection.
$ws.Range("F1").Select()

# Nudge page setup (orientation) for the now-wider sheet.
$ws.PageSetup.Orientation = 1   # xlPortrait
